$wb = $excel.ActiveWorkbook

# --- "Top Scores" sheet: add rows 5 and 6 ---
$ws1 = $wb.Worksheets.Item("Top Scores")

$ws1.Range("A5").Value = "Classic"
$ws1.Range("B5").Value = "Easy"
$ws1.Range("C5").Value = 141
$ws1.Range("D5").Value = 35

$ws1.Range("A6").Value = "Classic"
$ws1.Range("B6").Value = "Ultra Hard"
$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = 0

# Row 4 carries no explicit cell style (style 0); reuse it as the format
# donor for the new rows so the new cells don't inherit the column's style.
$ws1.Range("A4:D4").Copy()
$ws1.Range("A5:D6").PasteSpecial(-4122)

# --- "Classic" sheet: add rows 6, 7 and 8 ---
$ws2 = $wb.Worksheets.Item("Classic")

$ws2.Range("A6").Value = "Easy"
$ws2.Range("B6").Value = 141
$ws2.Range("C6").Value = 35

$ws2.Range("A7").Value = "Easy"
$ws2.Range("B7").Value = 211
$ws2.Range("C7").Value = 35

$ws2.Range("A8").Value = "Ultra Hard"
$ws2.Range("B8").Value = 1
$ws2.Range("C8").Value = 0

# Row 5 carries no explicit cell style (style 0); reuse it as the format
# donor for the new rows so the new cells don't inherit the column's style.
$ws2.Range("A5:C5").Copy()
$ws2.Range("A6:C8").PasteSpecial(-4122)
